$p = $ppt.ActivePresentation

# Insert a new "Title Only" slide right before the existing "Atom editor"
# slides (which sit at positions 29-30), matching the commit's
# "02. Hello World! add" insertion point.
$new = $p.Slides.Add(29, 7)

$tr = $new.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Code Editor"
$tr.Font.Italic = $true
